$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'42.524.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.92%  "

# Row 3
$ws.Range("D3").Value = "'2.232.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.00%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'114.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.41%  "

# Row 6
$ws.Range("D6").Value = "'281.75"
$ws.Range("D6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'0.628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.49%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("D9").Value = "'0.610"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.19%  "

# Row 10
$ws.Range("E10").Value = "  +0.29%  "

# Row 11
$ws.Range("D11").Value = "'0.0928"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.74%  "

# Row 12
$ws.Range("D12").Value = "'9.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.91%  "

# Row 13
$ws.Range("E13").Value = "  -2.93%  "

# Row 14
$ws.Range("D14").Value = "'15.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.55%  "

# Row 15
$ws.Range("D15").Value = "'0.880"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.19%  "

# Row 16
$ws.Range("D16").Value = "'2.569.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.88%  "

# Row 17
$ws.Range("D17").Value = "'2.241.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.02%  "

# Row 18
$ws.Range("D18").Value = "'42.709.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.16%  "

# Row 19
$ws.Range("E19").Value = "  -0.56%  "

# Row 20
$ws.Range("D20").Value = "'6.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.25%  "

# Row 21
$ws.Range("D21").Value = "'72.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.03%  "

# Row 22
$ws.Range("E22").Value = "  -3.81%  "

# Row 23
$ws.Range("D23").Value = "'3.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.21%  "

# Row 24
$ws.Range("D24").Value = "'231.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.21%  "

# Row 25
$ws.Range("E25").Value = "  -0.46%  "

# Row 26
$ws.Range("D26").Value = "'12.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.12%  "

# Row 27
$ws.Range("E27").Value = "  -1.66%  "

# Row 28
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").Value = "'40.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.58%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.06%  "

# Row 30
$ws.Range("B30").Value = "WEMIXToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D30").Value = "'3.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.19%  "

# Row 31
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "'173.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.07%  "

# Row 32
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'21.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.76%  "

# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0897"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.00%  "

# Row 34
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "'4.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +16.43%  "

# Row 35
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'5.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.60%  "

# Row 36
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").Value = "'0.128"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.31%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.0371"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.87%  "

# Row 38
$ws.Range("D38").Value = "'4.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.71%  "

# Row 39
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.106"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.64%  "

# Row 40
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'2.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.25%  "

# Row 41
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").Value = "'72.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.14%  "

# Row 42
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "'13.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.76%  "

# Row 43
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.234"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.63%  "

# Row 44
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.22%  "

# Row 45
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "'1.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.74%  "

# Row 46
$ws.Range("B46").Value = "THORChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D46").Value = "'5.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.18%  "

# Row 47
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").Value = "'1.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.45%  "

# Row 48
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'8.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.01%  "

# Row 49
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").Value = "'0.651"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.71%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0987"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.57%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'100.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.05%  "
